$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-10, columns A (Sending cluster) through T.
# For each sending cluster (ECs, FAPs, sCs) there are now 3 target-cluster
# rows (ECs, FAPs, sCs) instead of 2 -- NATMI recomputed with the new
# "sCs" cluster added to the analysis (per Dr Hou's advice).
$rows = @(
    @("ECs", "Adam12", "Sdc4", "ECs", 2, 0.6666666666666666, 1.958141333333333, 5.874423999999999, 0.1445807708852573, 0.1445807708852573, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 4.505505669852444, 40.54955102867199, 0.005460786284726423, 0.005460786284726423),
    @("ECs", "Adam12", "Sdc4", "FAPs", 2, 0.6666666666666666, 1.958141333333333, 5.874423999999999, 0.1445807708852573, 0.1445807708852573, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 33.21618967786222, 298.9457071007599, 0.04025885800953433, 0.04025885800953433),
    @("ECs", "Adam12", "Sdc4", "sCs", 2, 0.6666666666666666, 1.958141333333333, 5.874423999999999, 0.1445807708852573, 0.1445807708852573, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 81.56689223117066, 734.1020300805359, 0.09886112659099651, 0.09886112659099652),
    @("FAPs", "Adam12", "Sdc4", "ECs", 3, 1, 5.833003000000001, 17.499009, 0.4306839633891008, 0.4306839633891009, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 13.42121104406133, 120.790899396552, 0.01626684562495051, 0.01626684562495051),
    @("FAPs", "Adam12", "Sdc4", "FAPs", 3, 1, 5.833003000000001, 17.499009, 0.4306839633891008, 0.4306839633891009, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 98.94593957103167, 890.513456139285, 0.1199249694333544, 0.1199249694333544),
    @("FAPs", "Adam12", "Sdc4", "sCs", 3, 1, 5.833003000000001, 17.499009, 0.4306839633891008, 0.4306839633891009, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 242.975274044789, 2186.777466403101, 0.2944921483307959, 0.294492148330796),
    @("sCs", "Adam12", "Sdc4", "ECs", 3, 1, 5.752436333333333, 17.257309, 0.4247352657256419, 0.4247352657256419, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 13.23583444877244, 119.122510038952, 0.01604216452514934, 0.01604216452514934),
    @("sCs", "Adam12", "Sdc4", "FAPs", 3, 1, 5.752436333333333, 17.257309, 0.4247352657256419, 0.4247352657256419, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 97.57927740208721, 878.2134966187849, 0.1182685404828897, 0.1182685404828897),
    @("sCs", "Adam12", "Sdc4", "sCs", 3, 1, 5.752436333333333, 17.257309, 0.4247352657256419, 0.4247352657256419, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 239.6192483557556, 2156.573235201801, 0.2904245607176028, 0.2904245607176029)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $rowNum = $i + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}
